$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.45
$ws.Range("N2").Value = 3.5
$ws.Range("U2").Value = 1.58
$ws.Range("AB2").Value = 6.4

# Row 3
$ws.Range("L3").Value = 1.39
$ws.Range("N3").Value = 3.75

# Row 4
$ws.Range("F4").Value = 1.41
$ws.Range("L4").Value = 1.23
$ws.Range("N4").Value = 7
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 2.16
$ws.Range("U4").Value = 2.42
$ws.Range("AA4").Value = 230
$ws.Range("AD4").Value = 30
$ws.Range("AE4").Value = 95

# Row 5
$ws.Range("N5").Value = 4.6
$ws.Range("X5").Value = 18.5
$ws.Range("AA5").Value = 60

# Row 6
$ws.Range("G6").Value = 1.82
$ws.Range("P6").Value = 1.81
$ws.Range("W6").Value = 2.2
$ws.Range("AI6").Value = 95
$ws.Range("AL6").Value = 42
$ws.Range("AM6").Value = 150

# Row 7
$ws.Range("L7").Value = 1.41
$ws.Range("O7").Value = 1.32

# Row 8
$ws.Range("AN8").Value = 14.5

# Row 9
$ws.Range("S9").Value = 3.15
$ws.Range("T9").Value = 1.72
$ws.Range("X9").Value = 16

# Row 10
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.35
$ws.Range("L10").Value = 1.51
$ws.Range("X10").Value = 9.6
$ws.Range("AH10").Value = 20
$ws.Range("AO10").Value = 50

# Row 11
$ws.Range("L11").Value = 1.52
$ws.Range("U11").Value = 1.92

# Row 12
$ws.Range("L12").Value = 1.32
$ws.Range("P12").Value = 2.4
$ws.Range("AK12").Value = 140
$ws.Range("AM12").Value = 140

# Row 13
$ws.Range("I13").Value = 5.7
$ws.Range("P13").Value = 2.14
$ws.Range("U13").Value = 2.1

# Row 14
$ws.Range("H14").Value = 1.74
$ws.Range("I14").Value = 1.75
$ws.Range("N14").Value = 4.5
$ws.Range("X14").Value = 17
